$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark more Gantt cells as "Realizado" (done) / "Puesta en comun" (bordered milestone)
# for rows 10, 11, 12 (PYTHON - Vistas / Modelos / URLs), columns J:N
foreach ($row in 10, 11, 12) {
    $ws.Range("J$row").Value = "Realizado"
    [void]$ws.Range("C7").Copy()
    [void]$ws.Range("J$row").PasteSpecial(-4122)

    $ws.Range("K$row").Value = "Puesta en común"
    [void]$ws.Range("E7").Copy()
    [void]$ws.Range("K$row").PasteSpecial(-4122)

    $ws.Range("L$row").Value = "Realizado"
    [void]$ws.Range("C7").Copy()
    [void]$ws.Range("L$row").PasteSpecial(-4122)

    $ws.Range("M$row").Value = "Realizado"
    [void]$ws.Range("C7").Copy()
    [void]$ws.Range("M$row").PasteSpecial(-4122)

    $ws.Range("N$row").Value = "Puesta en común"
    [void]$ws.Range("E7").Copy()
    [void]$ws.Range("N$row").PasteSpecial(-4122)
}

# --- Row 13 (PYTHON - Bases de datos - Meter datos): update K/L and add new M/N/O markers
$ws.Range("K13").Value = "Realizado"
[void]$ws.Range("C7").Copy()
[void]$ws.Range("K13").PasteSpecial(-4122)

$ws.Range("L13").Value = "Puesta en común"
[void]$ws.Range("E7").Copy()
[void]$ws.Range("L13").PasteSpecial(-4122)

$ws.Range("M13").Value = "Retrasado"
[void]$ws.Range("C4").Copy()
[void]$ws.Range("M13").PasteSpecial(-4122)

$ws.Range("N13").Value = "Retrasado"
[void]$ws.Range("C4").Copy()
[void]$ws.Range("N13").PasteSpecial(-4122)

$ws.Range("O13").Value = "Pendiente"
[void]$ws.Range("C3").Copy()
[void]$ws.Range("O13").PasteSpecial(-4122)

# --- Remove the old "Cada uno que se encargue de una entidad de la entidad relación" task row
[void]$ws.Rows.Item(27).Delete()

# --- Move the "today" marker connector line further right (more tasks completed)
$shp = $ws.Shapes.Item(1)
$shp.Left = 1071
$shp.Top = 55.5

# --- Update the saved selection
[void]$ws.Range("G21").Select()
